$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "de olika algoritmerna " -> "dem " inside the paragraph that ends
#    with "Vad använder man de olika algoritmerna till, och hur
#    fungerar de?" (keeps the leading "d" and trailing "till, ..." in
#    their own runs, matching the target run layout).
# ------------------------------------------------------------------
$para = $null
$paraIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*de olika algoritmerna*") {
        $para = $p
        $paraIndex = $i
        break
    }
}

$fullText = $para.Range.Text
$marker = "e olika algoritmerna "
$offset = $fullText.IndexOf($marker)
$segStart = $para.Range.Start + $offset
$segEnd = $segStart + $marker.Length

# Replace "e olika algoritmerna " with "em " in place.
$seg = $d.Range($segStart, $segEnd)
$seg.Text = "em "

# Force the freshly typed "em " text to live in its own run (no
# formatting differences) by briefly bracketing it with a throwaway
# bookmark, then deleting that bookmark again. Deleting a bookmark
# only removes the bookmark markers, not the run split it produced.
$newSegEnd = $segStart + 3
$tempRange = $d.Range($segStart, $newSegEnd)
$d.Bookmarks.Add("__tmp_split__", $tempRange)
$d.Bookmarks("__tmp_split__").Delete()

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the last (empty) paragraph to
#    the empty paragraph immediately following the edited paragraph.
#    Re-adding a bookmark with the same name simply relocates it, so
#    it disappears from its old spot automatically.
# ------------------------------------------------------------------
$bookmarkHome = $d.Paragraphs.Item($paraIndex + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkHome.Range)
